$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.746.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.22%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'1.861.36"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -0.70%  "
$ws.Range("E3").ClearFormats()

$ws.Range("D4").Value = "'1.014"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +1.02%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = "'333.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.56%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D7").Value = "'0.4701"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -0.29%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").Value = "'0.3893"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -1.25%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").Value = "'46.66"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -2.66%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'0.07968"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -0.76%  "
$ws.Range("E10").ClearFormats()

$ws.Range("E11").Value = "'  -2.31%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'21.57"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -2.08%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'1.874.13"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +1.19%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'5.985"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +0.42%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = "'7.123"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.10%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = "'1.015"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.84%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = "'88.20"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.38%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'0.06682"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.08%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = "'0.00001042"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -0.61%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").Value = "'16.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -1.50%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").Value = "'1.011"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +0.78%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = "'27.754.65"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.25%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'5.462"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -0.97%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'10.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.51%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = "'2.323"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +0.73%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").Value = "'2.097.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +1.07%  "
$ws.Range("E26").ClearFormats()

$ws.Range("D27").Value = "'158.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -0.06%  "
$ws.Range("E27").ClearFormats()

$ws.Range("D28").Value = "'19.64"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -2.52%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").Value = "'2.086"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -0.75%  "
$ws.Range("E29").ClearFormats()

$ws.Range("D30").Value = "'5.404"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -2.95%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").Value = "'120.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -1.04%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").Value = "'0.9653"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -0.87%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").Value = "'0.09439"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -0.88%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = "'3.640"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +1.26%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = "'5.298"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -0.55%  "
$ws.Range("E35").ClearFormats()

$ws.Range("E36").Value = "'  -7.26%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").Value = "'0.06034"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -1.03%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").Value = "'0.02219"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -1.37%  "
$ws.Range("E38").ClearFormats()

$ws.Range("E39").Value = "'  -1.80%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'8.128"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -0.79%  "
$ws.Range("E40").ClearFormats()

$ws.Range("D41").Value = "'1.011"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +0.85%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'0.5916"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -1.76%  "
$ws.Range("E42").ClearFormats()

$ws.Range("D43").Value = "'0.1882"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -1.09%  "
$ws.Range("E43").ClearFormats()

$ws.Range("D44").Value = "'10.20"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -0.19%  "
$ws.Range("E44").ClearFormats()

$ws.Range("D45").Value = "'1.260"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -0.65%  "
$ws.Range("E45").ClearFormats()

$ws.Range("D46").Value = "'0.5618"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -1.26%  "
$ws.Range("E46").ClearFormats()

$ws.Range("D47").Value = "'12.05"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.96%  "
$ws.Range("E47").ClearFormats()

$ws.Range("D48").Value = "'1.916"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -1.41%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").Value = "'3.305"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -2.20%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").Value = "'0.06770"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -1.67%  "
$ws.Range("E50").ClearFormats()

$ws.Range("D51").Value = "'112.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -2.52%  "
$ws.Range("E51").ClearFormats()
